# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Datos actualizados" timestamp (A1): 21:05 -> 21:35 ---
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 21:35"

# --- Estados Unidos (row 4): updated case counts ---
$ws.Range("B4").Value = 1396184
$ws.Range("C4").Value = 10350
$ws.Range("E4").Value = 1038491
$ws.Range("G4").Value = 996
$ws.Range("H4").Value = 82791

# --- Dinamarca (row 46): updated case counts ---
$ws.Range("D46").Value = 8580
$ws.Range("E46").Value = 1484
$ws.Range("G46").Value = 0
$ws.Range("H46").Value = 527

# --- Costa Rica overtakes Gabon (rows 108/109 swap order + new data) ---
$ws.Range("A108").Value = "Costa Rica"
$ws.Range("B108").Value = 804
$ws.Range("C108").Value = 3
$ws.Range("D108").Value = 520
$ws.Range("E108").Value = 277
$ws.Range("F108").Value = 6
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 7

$ws.Range("A109").Value = "Gabon"
$ws.Range("B109").Value = 802
$ws.Range("C109").Value = 0
$ws.Range("D109").Value = 127
$ws.Range("E109").Value = 666
$ws.Range("F109").Value = 1
$ws.Range("G109").Value = 0
$ws.Range("H109").Value = 9

# --- Curazao overtakes Dominica (rows 198/199 swap order) ---
$ws.Range("A198").Value = "Curazao"
$ws.Range("B198").Value = 16
$ws.Range("C198").Value = 0
$ws.Range("D198").Value = 14
$ws.Range("E198").Value = 1
$ws.Range("F198").Value = 0
$ws.Range("G198").Value = 0
$ws.Range("H198").Value = 1

$ws.Range("A199").Value = "Dominica"
$ws.Range("B199").Value = 16
$ws.Range("C199").Value = 0
$ws.Range("D199").Value = 15
$ws.Range("E199").Value = 1
$ws.Range("F199").Value = 0
$ws.Range("G199").Value = 0
$ws.Range("H199").Value = 0
